$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-102 down to 28-103
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new record
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C27").Value = "Ñuble"
$ws.Range("D27").Value = "2023-07-26"
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112001
$ws.Range("G27").Value = "Berenjena"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 9000
$ws.Range("N27").Value = "$/caja 60 unidades"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 150
$ws.Range("Q27").Value = 60
$ws.Range("R27").Value = "Hortaliza"
